$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "70.495.90"
Set-TextValue $ws.Range("E2") "  +0.21%  "
Set-TextValue $ws.Range("D3") "3.613.26"
Set-TextValue $ws.Range("E3") "  -0.52%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "581.27"
Set-TextValue $ws.Range("E5") "  -1.73%  "
Set-TextValue $ws.Range("D6") "190.46"
Set-TextValue $ws.Range("E6") "  -2.11%  "
Set-TextValue $ws.Range("D7") "3.609.04"
Set-TextValue $ws.Range("E7") "  -0.47%  "
Set-TextValue $ws.Range("D8") "0.631"
Set-TextValue $ws.Range("E8") "  -2.07%  "
Set-TextValue $ws.Range("E9") "  +0.07%  "
Set-TextValue $ws.Range("D10") "0.189"
Set-TextValue $ws.Range("E10") "  +4.77%  "
Set-TextValue $ws.Range("E11") "  -1.31%  "
Set-TextValue $ws.Range("D12") "56.16"
Set-TextValue $ws.Range("E12") "  -4.04%  "
Set-TextValue $ws.Range("D13") "0.0000313"
Set-TextValue $ws.Range("E13") "  +7.73%  "
Set-TextValue $ws.Range("D14") "9.71"
Set-TextValue $ws.Range("E14") "  -2.33%  "
Set-TextValue $ws.Range("D15") "4.193.60"
Set-TextValue $ws.Range("E15") "  -0.46%  "
Set-TextValue $ws.Range("D16") "19.83"
Set-TextValue $ws.Range("E16") "  -0.28%  "
Set-TextValue $ws.Range("D17") "3.609.60"
Set-TextValue $ws.Range("E17") "  -0.67%  "
Set-TextValue $ws.Range("D18") "70.458.44"
Set-TextValue $ws.Range("E18") "  +0.23%  "
Set-TextValue $ws.Range("D19") "12.73"
Set-TextValue $ws.Range("E19") "  +0.15%  "
Set-TextValue $ws.Range("E20") "  +0.15%  "
Set-TextValue $ws.Range("D21") "1.05"
Set-TextValue $ws.Range("E21") "  -1.87%  "
Set-TextValue $ws.Range("D22") "493.02"
Set-TextValue $ws.Range("E22") "  +0.96%  "
Set-TextValue $ws.Range("D23") "19.37"
Set-TextValue $ws.Range("E23") "  +0.38%  "
Set-TextValue $ws.Range("E24") "  -8.16%  "
Set-TextValue $ws.Range("D25") "96.38"
Set-TextValue $ws.Range("E25") "  +5.40%  "
Set-TextValue $ws.Range("D26") "4.38"
Set-TextValue $ws.Range("E26") "  -2.05%  "
Set-TextValue $ws.Range("D27") "2.99"
Set-TextValue $ws.Range("E27") "  -5.58%  "
Set-TextValue $ws.Range("D28") "11.11"
Set-TextValue $ws.Range("E28") "  -3.58%  "
Set-TextValue $ws.Range("D29") "9.46"
Set-TextValue $ws.Range("E29") "  -1.52%  "
Set-TextValue $ws.Range("D30") "32.30"
Set-TextValue $ws.Range("E30") "  -2.01%  "
Set-TextValue $ws.Range("D31") "7.63"
Set-TextValue $ws.Range("E31") "  -3.87%  "
Set-TextValue $ws.Range("D32") "12.23"
Set-TextValue $ws.Range("E32") "  -0.24%  "
Set-TextValue $ws.Range("B33") "Hedera"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D33") "0.118"
Set-TextValue $ws.Range("E33") "  -3.07%  "
Set-TextValue $ws.Range("B34") "OKB"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D34") "65.93"
Set-TextValue $ws.Range("E34") "  +0.22%  "
Set-TextValue $ws.Range("D35") "574.43"
Set-TextValue $ws.Range("E35") "  -8.30%  "
Set-TextValue $ws.Range("D36") "38.49"
Set-TextValue $ws.Range("E36") "  -6.31%  "
Set-TextValue $ws.Range("D37") "0.0₃0818"
Set-TextValue $ws.Range("E37") "  -0.42%  "
Set-TextValue $ws.Range("E38") "  +0.20%  "
Set-TextValue $ws.Range("D39") "3.37"
Set-TextValue $ws.Range("E39") "  +18.15%  "
Set-TextValue $ws.Range("D40") "0.398"
Set-TextValue $ws.Range("E40") "  -3.51%  "
Set-TextValue $ws.Range("D41") "3.00"
Set-TextValue $ws.Range("E41") "  +5.54%  "
Set-TextValue $ws.Range("D42") "3.55"
Set-TextValue $ws.Range("E42") "  -0.64%  "
Set-TextValue $ws.Range("E43") "  -6.05%  "
Set-TextValue $ws.Range("D44") "3.04"
Set-TextValue $ws.Range("E44") "  -4.40%  "
Set-TextValue $ws.Range("D45") "3.57"
Set-TextValue $ws.Range("E45") "  +8.23%  "
Set-TextValue $ws.Range("D46") "3.233.59"
Set-TextValue $ws.Range("E46") "  -1.90%  "
Set-TextValue $ws.Range("D47") "0.0445"
Set-TextValue $ws.Range("E47") "  -2.13%  "
Set-TextValue $ws.Range("D48") "9.82"
Set-TextValue $ws.Range("E48") "  +6.94%  "
Set-TextValue $ws.Range("E49") "  -0.32%  "
Set-TextValue $ws.Range("E50") "  -0.11%  "
Set-TextValue $ws.Range("D51") "3.21"
Set-TextValue $ws.Range("E51") "  -3.85%  "
